# Ticket 82 - "NoSpaceAfterParen" test sheet: update cell reference in a
# formula even when no space precedes it, e.g. $[A2-(IF(B2="-",0,B2)+C2)].
#
# Add a new worksheet at the end of the workbook (after the last existing
# sheet) named "NoSpaceAfterParen" and populate it with the header row and
# the jt:for / formula test data used to exercise the fix.

$wb = $excel.ActiveWorkbook

$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws.Name = "NoSpaceAfterParen"

# Header row - bold, centered.
$ws.Range("A1").Value = "a"
$ws.Range("B1").Value = "b"
$ws.Range("C1").Value = "c"
$ws.Range("D1").Value = "result"
$ws.Range("A1:D1").Font.Bold = $true
$ws.Range("A1:D1").HorizontalAlignment = -4108  # xlCenter

# Data row - jt:for loop plus formula tags.
$ws.Range("A2").Value = '<jt:for start="1" end="10" var="x">${x}'
$ws.Range("B2").Value = '${x+1}'
$ws.Range("C2").Value = '${x+2}'
$ws.Range("D2").Value = '$[A2-(IF(B2="-",0,B2)+C2)]'
$ws.Range("E2").Value = '</jt:for>'

# Restore the selection on a couple of other sheets that were left with a
# different selection the last time the workbook was saved.
$wsMultiLevel2 = $wb.Worksheets.Item("MultiLevel2")
$wsMultiLevel2.Range("E6").Select() | Out-Null

$wsTagParseInFormula = $wb.Worksheets.Item("TagParseInFormula")
$wsTagParseInFormula.Range("A3:XFD3").Select() | Out-Null

# Leave the first sheet active/selected, as in the original workbook.
$wb.Worksheets.Item("Formula Test").Select() | Out-Null
